$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) matching the style already used by the other
# header cells (e.g. H1 - bold font, thin border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (row 2), plain numeric values like the rest of row 2.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
